$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 32: USerMaster section header (no border style) ---
$ws.Range("A32").Value = "USerMaster"

# --- Row 31: update the BranchMaster sample-data row ---
$ws.Range("B31").Value = "TestBranch"
$ws.Range("C31").Value = "TestAdd"
# D31 "Address line 2", E31 "TestCity", J31 "Contact Name" stay textually the same
$ws.Range("F31").Value = "TestState"
$ws.Range("G31").Value = "TestZip"
$ws.Range("H31").Value = "TestCountry"
$ws.Range("K31").Value = 9999999999
$ws.Range("L31").Value = "Test@ctdi.com"

# --- New row 34: BinsMaster section header (no border style) ---
$ws.Range("A34").Value = "BinsMaster"

# --- New row 33: UserMaster sample-data row (bordered, like row 31) ---
$ws.Range("A33").Value = "Test"
$ws.Range("B33").Value = "User"
$ws.Range("C33").Value = "Tuser"
$ws.Range("D33").Value = "Tuser@gmail.com"
$ws.Range("E33").Value = "Developer"
$ws.Range("F33").Value = "Test"
$ws.Range("G33").Value = "ViewOrders"

# Copy the bordered style from the row-31 data cells onto the new row-33 cells
$ws.Range("B31").Copy()
$ws.Range("A33:G33").PasteSpecial(-4122)

# --- Hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("L31"), "mailto:Test@ctdi.com")
$ws.Range("K31").Copy()
$ws.Range("L31").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("D33"), "mailto:Tuser@gmail.com")
$ws.Range("C33").Copy()
$ws.Range("D33").PasteSpecial(-4122)

# --- Selection state ---
$ws.Range("A33:G33").Select()

Write-Output "edit applied"
